# Applies the "Description.docx" update:
#   - Paragraph 1 ("Rgwerh ryrtbyrtvt") -> "Disclaimer" (bold heading)
#   - New paragraph inserted after it with the disclaimer body text
#   - The old "T" paragraph becomes empty, and two more empty paragraphs
#     are added so the blank-line gap grows from 1 to 6 paragraphs
#   - Final heading "Description" -> "How to use"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Turn paragraph 1 into the bold "Disclaimer" heading and add a new
#    paragraph right after it for the disclaimer body copy.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphBefore()
# Paragraph layout is now:
#   1 = brand-new empty paragraph (will become "Disclaimer")
#   2 = original "Rgwerh ryrtbyrtvt" paragraph (with its proofErr runs)

$headingPara = $d.Paragraphs.Item(1)
$headingPara.Range.InsertParagraphAfter()
# Paragraph layout is now:
#   1 = brand-new empty paragraph (will become "Disclaimer")
#   2 = brand-new empty paragraph (will become the body text)
#   3 = original "Rgwerh ryrtbyrtvt" paragraph

$bodyPara = $d.Paragraphs.Item(2)
$bodyPara.Range.Text = "This tool was designed to get a sense of the differences in electricity consumption between different products. It’s often difficult to understand whether activities matter a lot or very little for our overall energy consumption."

# Drop the old "Rgwerh ryrtbyrtvt" paragraph entirely (text + proofErr markers).
$oldPara = $d.Paragraphs.Item(3)
$oldPara.Range.Delete()

# Fill in the heading text and bold it.
$headingPara = $d.Paragraphs.Item(1)
$headingPara.Range.Text = "Disclaimer"
$headingPara.Range.Font.Bold = $true
$headingPara.Range.Font.BoldBi = $true

# ---------------------------------------------------------------------
# 2) The paragraph that used to hold "T" becomes blank, and two more
#    blank paragraphs are inserted after it.
# ---------------------------------------------------------------------
$tPara = $d.Paragraphs.Item(4)
$textOnly = $d.Range($tPara.Range.Start, $tPara.Range.End - 1)
$textOnly.Delete()

$tPara = $d.Paragraphs.Item(4)
$tPara.Range.InsertParagraphAfter()
$extraBlank = $d.Paragraphs.Item(5)
$extraBlank.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 3) Rename the closing "Description" heading to "How to use" (a plain
#    text swap keeps its existing bold run/paragraph formatting intact).
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute("Description", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "How to use", 2)
